$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Contract/order number (2 identical occurrences)
Replace-AllText "2487097" "1234568"

# Contract/order date (2 identical occurrences)
Replace-AllText "03.05.2024" "09.02.2024"

# Vessel name "Светлый" -> "Светлый " (2 identical occurrences, trailing space added)
Replace-AllText "Светлый" "Светлый "

# Signatory / title block
Replace-AllText "генеральный директор Котлярчук О. Е." "Капитан Бахтин Ю. Г."

# Basis of authority
Replace-AllText "Устава" "Кодекса торгового мореплавания (КТМ РФ) "

# Survey description
Replace-AllText "Освидетельствование объектов технаблюдения:Насосы Q=64 m3/ч - 2 шт.Насосы Q=36 м3/ч - 3 шт.Компенсаторы DN-200 - 15 шт." "Освидетельствование // Survey of blah blah,  blah blah blah,  blah blah blah blah"

# Certificate numbers / date
Replace-AllText "Свидетельство ф. 6.5.30 №№ 24.42.03.00234.121 - 24.42.03.00236.121 от 06.05.2024" "Свидетельство ф. 6.5.30 №№ 24.42.03.00414.121 - 24.42.03.00416.121 от --"

# Cost of services
Replace-AllText "10 000,00 p. (десять тысяч рублей 00 копеек)" "536 112,20 p. (пятьсот тридцать шесть тысяч сто двенадцать рублей 20 копеек)"

# Final signature block
Replace-AllText "О. Е. Котлярчук" "Ю. Г. Бахтин"
